$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.139906333333333
$ws.Range("H2").Value = 3.419719
$ws.Range("I2").Value = 0.2178538649973528
$ws.Range("J2").Value = 0.2178538649973527
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 101.7913436666667
$ws.Range("N2").Value = 305.374031
$ws.Range("O2").Value = 0.2805454758424659
$ws.Range("P2").Value = 0.2805454758424659
$ws.Range("Q2").Value = 116.0325973241432
$ws.Range("R2").Value = 1044.293375917289
$ws.Range("S2").Value = 0.06111791621980266
$ws.Range("T2").Value = 0.06111791621980264

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.139906333333333
$ws.Range("H3").Value = 3.419719
$ws.Range("I3").Value = 0.2178538649973528
$ws.Range("J3").Value = 0.2178538649973527
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 54.34621166666667
$ws.Range("N3").Value = 163.038635
$ws.Range("O3").Value = 0.1497827149446808
$ws.Range("P3").Value = 0.1497827149446808
$ws.Range("Q3").Value = 61.94959087150722
$ws.Range("R3").Value = 557.5463178435649
$ws.Range("S3").Value = 0.03263074336049546
$ws.Range("T3").Value = 0.03263074336049546

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.139906333333333
$ws.Range("H4").Value = 3.419719
$ws.Range("I4").Value = 0.2178538649973528
$ws.Range("J4").Value = 0.2178538649973527
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 122.2542826666667
$ws.Range("N4").Value = 366.762848
$ws.Range("O4").Value = 0.3369430510399163
$ws.Range("P4").Value = 0.3369430510399163
$ws.Range("Q4").Value = 139.3584310888569
$ws.Range("R4").Value = 1254.225879799712
$ws.Range("S4").Value = 0.07340434595304608
$ws.Range("T4").Value = 0.07340434595304605

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.139906333333333
$ws.Range("H5").Value = 3.419719
$ws.Range("I5").Value = 0.2178538649973528
$ws.Range("J5").Value = 0.2178538649973527
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 75.13712199999999
$ws.Range("N5").Value = 225.411366
$ws.Range("O5").Value = 0.2070842066291165
$ws.Range("P5").Value = 0.2070842066291166
$ws.Range("Q5").Value = 85.64928123623932
$ws.Range("R5").Value = 770.8435311261539
$ws.Range("S5").Value = 0.04511409479406345
$ws.Range("T5").Value = 0.04511409479406345

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.139906333333333
$ws.Range("H6").Value = 3.419719
$ws.Range("I6").Value = 0.2178538649973528
$ws.Range("J6").Value = 0.2178538649973527
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.304706666666666
$ws.Range("N6").Value = 27.91412
$ws.Range("O6").Value = 0.02564455154382035
$ws.Range("P6").Value = 0.02564455154382035
$ws.Range("Q6").Value = 10.60649405914222
$ws.Range("R6").Value = 95.45844653227998
$ws.Range("S6").Value = 0.005586764669945094
$ws.Range("T6").Value = 0.005586764669945092

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.902924
$ws.Range("H7").Value = 5.708772
$ws.Range("I7").Value = 0.3636784322304457
$ws.Range("J7").Value = 0.3636784322304456
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 101.7913436666667
$ws.Range("N7").Value = 305.374031
$ws.Range("O7").Value = 0.2805454758424659
$ws.Range("P7").Value = 0.2805454758424659
$ws.Range("Q7").Value = 193.701190855548
$ws.Range("R7").Value = 1743.310717699932
$ws.Range("S7").Value = 0.1020283388237324
$ws.Range("T7").Value = 0.1020283388237324

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.902924
$ws.Range("H8").Value = 5.708772
$ws.Range("I8").Value = 0.3636784322304457
$ws.Range("J8").Value = 0.3636784322304456
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.34621166666667
$ws.Range("N8").Value = 163.038635
$ws.Range("O8").Value = 0.1497827149446808
$ws.Range("P8").Value = 0.1497827149446808
$ws.Range("Q8").Value = 103.41671048958
$ws.Range("R8").Value = 930.7503944062199
$ws.Range("S8").Value = 0.05447274294630126
$ws.Range("T8").Value = 0.05447274294630126

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.902924
$ws.Range("H9").Value = 5.708772
$ws.Range("I9").Value = 0.3636784322304457
$ws.Range("J9").Value = 0.3636784322304456
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 122.2542826666667
$ws.Range("N9").Value = 366.762848
$ws.Range("O9").Value = 0.3369430510399163
$ws.Range("P9").Value = 0.3369430510399163
$ws.Range("Q9").Value = 232.640608589184
$ws.Range("R9").Value = 2093.765477302656
$ws.Range("S9").Value = 0.1225389205531398
$ws.Range("T9").Value = 0.1225389205531398

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.902924
$ws.Range("H10").Value = 5.708772
$ws.Range("I10").Value = 0.3636784322304457
$ws.Range("J10").Value = 0.3636784322304456
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 75.13712199999999
$ws.Range("N10").Value = 225.411366
$ws.Range("O10").Value = 0.2070842066291165
$ws.Range("P10").Value = 0.2070842066291166
$ws.Range("Q10").Value = 142.980232744728
$ws.Range("R10").Value = 1286.822094702552
$ws.Range("S10").Value = 0.07531205960656276
$ws.Range("T10").Value = 0.07531205960656276

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.902924
$ws.Range("H11").Value = 5.708772
$ws.Range("I11").Value = 0.3636784322304457
$ws.Range("J11").Value = 0.3636784322304456
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.304706666666666
$ws.Range("N11").Value = 27.91412
$ws.Range("O11").Value = 0.02564455154382035
$ws.Range("P11").Value = 0.02564455154382035
$ws.Range("Q11").Value = 17.70614962896
$ws.Range("R11").Value = 159.35534666064
$ws.Range("S11").Value = 0.009326370300709442
$ws.Range("T11").Value = 0.009326370300709439

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.189605333333333
$ws.Range("H12").Value = 6.568816
$ws.Range("I12").Value = 0.4184677027722017
$ws.Range("J12").Value = 0.4184677027722016
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 101.7913436666667
$ws.Range("N12").Value = 305.374031
$ws.Range("O12").Value = 0.2805454758424659
$ws.Range("P12").Value = 0.2805454758424659
$ws.Range("Q12").Value = 222.8828689796995
$ws.Range("R12").Value = 2005.945820817296
$ws.Range("S12").Value = 0.1173992207989309
$ws.Range("T12").Value = 0.1173992207989309

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.189605333333333
$ws.Range("H13").Value = 6.568816
$ws.Range("I13").Value = 0.4184677027722017
$ws.Range("J13").Value = 0.4184677027722016
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 54.34621166666667
$ws.Range("N13").Value = 163.038635
$ws.Range("O13").Value = 0.1497827149446808
$ws.Range("P13").Value = 0.1497827149446808
$ws.Range("Q13").Value = 118.9967549117955
$ws.Range("R13").Value = 1070.97079420616
$ws.Range("S13").Value = 0.06267922863788411
$ws.Range("T13").Value = 0.0626792286378841

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.189605333333333
$ws.Range("H14").Value = 6.568816
$ws.Range("I14").Value = 0.4184677027722017
$ws.Range("J14").Value = 0.4184677027722016
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 122.2542826666667
$ws.Range("N14").Value = 366.762848
$ws.Range("O14").Value = 0.3369430510399163
$ws.Range("P14").Value = 0.3369430510399163
$ws.Range("Q14").Value = 267.6886293497742
$ws.Range("R14").Value = 2409.197664147968
$ws.Range("S14").Value = 0.1409997845337305
$ws.Range("T14").Value = 0.1409997845337305

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.189605333333333
$ws.Range("H15").Value = 6.568816
$ws.Range("I15").Value = 0.4184677027722017
$ws.Range("J15").Value = 0.4184677027722016
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 75.13712199999999
$ws.Range("N15").Value = 225.411366
$ws.Range("O15").Value = 0.2070842066291165
$ws.Range("P15").Value = 0.2070842066291166
$ws.Range("Q15").Value = 164.5206430625173
$ws.Range("R15").Value = 1480.685787562656
$ws.Range("S15").Value = 0.08665805222849034
$ws.Range("T15").Value = 0.08665805222849034

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.189605333333333
$ws.Range("H16").Value = 6.568816
$ws.Range("I16").Value = 0.4184677027722017
$ws.Range("J16").Value = 0.4184677027722016
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 9.304706666666666
$ws.Range("N16").Value = 27.91412
$ws.Range("O16").Value = 0.02564455154382035
$ws.Range("P16").Value = 0.02564455154382035
$ws.Range("Q16").Value = 20.37363534243555
$ws.Range("R16").Value = 183.36271808192
$ws.Range("S16").Value = 0.01073141657316582
$ws.Range("T16").Value = 0.01073141657316582
